# Database input data fixed
# Append repeated/extra sample rows to the "Bugs" and "Sales" sheets, mirroring
# the existing rows so the demo data has more entries to work with.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $value) {
    # Force the cell to be stored as text even when the value looks numeric
    # (e.g. "10", "4", "69000"), matching the existing rows in this sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---- Sheet "Bugs" (sheet1) ----
$bugs = $wb.Worksheets.Item("Bugs")

$bugsRows = @(
    @("Opel Astra", "Break fluid leak", "10"),
    @("Opel Corsa", "Engine problem", "10"),
    @("Mitsubishi Space Star", "Thick-Thick Steering Wheel Sounds", "4"),
    @("Opel Astra", "Break fluid leak", "10"),
    @("Opel Corsa", "Engine problem", "10"),
    @("Mitsubishi Space Star", "Thick-Thick Steering Wheel Sounds", "4"),
    @("Opel Astra", "Break fluid leak", "10"),
    @("Opel Corsa", "Engine problem", "10"),
    @("Mitsubishi Space Star", "Thick-Thick Steering Wheel Sounds", "4")
)

$startRow = 5
for ($i = 0; $i -lt $bugsRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $bugsRows[$i]
    $bugs.Cells.Item($r, 1).Value = $rowData[0]
    $bugs.Cells.Item($r, 2).Value = $rowData[1]
    Set-TextCell $bugs.Cells.Item($r, 3) $rowData[2]
}

# ---- Sheet "Sales" (sheet2) ----
$sales = $wb.Worksheets.Item("Sales")

$salesRows = @(
    @("Astra", "Opel", "3", "69000"),
    @("Astra", "Opel", "3", "69000"),
    @("Astra", "Opel", "3", "69000"),
    @("Astra", "Opel", "4", "92000"),
    @("Astra", "Opel", "4", "92000"),
    @("Astra", "Opel", "4", "92000"),
    @("Corsa", "Opel", "2", "46000"),
    @("Corsa", "Opel", "2", "46000"),
    @("Corsa", "Opel", "2", "46000")
)

$startRow = 5
for ($i = 0; $i -lt $salesRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $salesRows[$i]
    $sales.Cells.Item($r, 1).Value = $rowData[0]
    $sales.Cells.Item($r, 2).Value = $rowData[1]
    Set-TextCell $sales.Cells.Item($r, 3) $rowData[2]
    Set-TextCell $sales.Cells.Item($r, 4) $rowData[3]
}
